# Auto commit at 2025-10-12  7:37:20.44
#
# 1. Update the "Metrics" sheet figures (B2:B13) with the new numbers.
# 2. Clear today's running totals (B3:B6) on the "today" sheet — they are
#    reset to blank (the day rolled over), which also makes the F-column
#    "today" formulas collapse onto the E-column totals since Bx now reads
#    as 0.
# 3. Move the active sheet / selection from "Metrics" (D7) to "today" (D4),
#    which is reflected in workbookView.activeTab + the per-sheet
#    tabSelected/selection state.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- 1. Refresh the Metrics figures ------------------------------------
$metrics.Range("B2").Value  = 156484.72
$metrics.Range("B3").Value  = 131240.27000000002
$metrics.Range("B4").Value  = 57031.14
$metrics.Range("B5").Value  = 6123
$metrics.Range("B6").Value  = 4523616.1900000004
$metrics.Range("B7").Value  = 3821058.9399999995
$metrics.Range("B8").Value  = 1327633.2800000003
$metrics.Range("B9").Value  = 175124
$metrics.Range("B10").Value = 32988939.990999825
$metrics.Range("B11").Value = 31096280.460000005
$metrics.Range("B12").Value = 11609342.169999998
$metrics.Range("B13").Value = 1272751

# --- 2. Clear today's running-total inputs ------------------------------
$today.Range("B3:B6").ClearContents() | Out-Null

# --- 3. Switch the active sheet / selection -----------------------------
$metrics.Activate() | Out-Null
$metrics.Range("D9").Select() | Out-Null

$today.Activate() | Out-Null
$today.Range("D4").Select() | Out-Null
